$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data for rows 2-11 (A: employee_id, B: employee_name, C: department,
# D: absence_reason, E: absence_duration, F: absence_date, G: salary)
$data = @(
    @(2,  82323, "Isabella Azevedo",      "P&D",                     "Viagem de negocios",   1, 45086,  7836.48),
    @(3,  31897, "Erick Cunha",           "Recursos Humanos",        "Outros",                1, 45078,  6141.26),
    @(4,  68578, "Ana Beatriz Moura",     "Recursos Humanos",        "Outros",                6, 45099,  7105.51),
    @(5,  96305, "Vitória Cirino",        "Engenharia",              "Consulta medica",       3, 45103,  3524.75),
    @(6,  96533, "Vicente Melo",          "Atendimento ao Cliente",  "Consulta medica",       6, 45098,  8722.059999999999),
    @(7,  87004, "Ana Vitória Cassiano",  "P&D",                     "Doenca",                4, 45105,  3816.99),
    @(8,  29362, "Antony Guerra",         "Engenharia",              "Outros",                4, 45086,  4933.94),
    @(9,  77140, "Luísa Vargas",          "Marketing",               "Consulta medica",       8, 45087,  3941.87),
    @(10, 46635, "Maria Helena da Cunha", "P&D",                     "Problemas pessoais",    6, 45100,  6780.21),
    @(11, 76859, "Marina Rodrigues",      "Juridico",                "Doenca",                1, 45087,  9404.719999999999)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
}
